$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) INTENT sheet: add a new intent row for "AMAZON.YesIntent" (row 10)
# ---------------------------------------------------------------------------
$intent = $wb.Worksheets.Item("INTENT")
$intent.Range("A10").Value = "AMAZON.YesIntent"

# ---------------------------------------------------------------------------
# 2) UTTERANCES_MAIN sheet: add a new column F for the "AMAZON.YesIntent"
#    utterances ("true" and "yes"), matching the header/value styling already
#    used by the other intent columns on this sheet.
# ---------------------------------------------------------------------------
$utt = $wb.Worksheets.Item("UTTERANCES_MAIN")

# Copy the header formatting from the neighbouring header cell (E1) onto the
# new header cell (F1), then set its text.
$utt.Range("E1").Copy()
$utt.Range("F1").PasteSpecial(-4122)
$utt.Range("F1").Value = "AMAZON.YesIntent"

# F2 holds the literal text "true" (not the boolean TRUE) - force a Text
# number format first and prefix the value with an apostrophe so it is
# stored as a shared string rather than being auto-coerced to a boolean.
$utt.Range("F2").NumberFormat = "@"
$utt.Range("F2").Value = "'true"

# F3 holds the plain utterance text "yes".
$utt.Range("F3").Value = "yes"

# ---------------------------------------------------------------------------
# 3) Reflect that UTTERANCES_MAIN was the sheet being worked on: make it the
#    active sheet/tab and leave a sensible selection on it.
# ---------------------------------------------------------------------------
$utt.Activate()
$utt.Range("E1").Select()
$utt.Range("F4").Select()
